$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bookFlight")

# --- Populate new rows 8-32 (write order chosen to reproduce original shared-string table order) ---
$ws.Range("A9").Value = "Last Name"
$ws.Range("A8").Value = "First Name"
$ws.Range("A10").Value = "Meal "
$ws.Range("B10").Value = "Vegetarian"
$ws.Range("A11").Value = "Card Type"
$ws.Range("B11").Value = "MasterCard"
$ws.Range("A12").Value = "Card Number"
$ws.Range("A13").Value = "Expiry Month "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "08"
$ws.Range("A14").Value = "Expiry Year"
$ws.Range("A16").Value = "Middle Name"
$ws.Range("B8").Value = "First"
$ws.Range("B9").Value = "Last"
$ws.Range("B16").Value = "Middle"
$ws.Range("A18").Value = "Ticketless Travel"
$ws.Range("B18").Value = "Yes"
$ws.Range("A19").Value = "Billing Address 1"
$ws.Range("A20").Value = "Billing Address 2"
$ws.Range("B21").Value = "Auckland"
$ws.Range("A23").Value = "Postal Code"
$ws.Range("A25").Value = "Delivery same as billing"
$ws.Range("A26").Value = "Delivery Address1"
$ws.Range("A27").Value = "Delivery Addrss2"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "1234 5678 9874 1523"
$ws.Range("B24").Value = "NEW ZEALAND"
$ws.Range("A32").Value = "Warning"
$ws.Range("B32").Value = "You have chosen a mailing location outside of the United States and its territories. An additional charge of `$6.5 will be added as mailing charge.`nYou have chosen a mailing location outside of the United States and its territories. An additional charge of `$6.5 will be added as mailing charge.`n"
$ws.Range("B32").WrapText = $true
$ws.Rows(32).RowHeight = 14.25
$ws.Range("B14").Value = 2000
$ws.Range("A15").Value = "First Name"
$ws.Range("B15").Value = "First"
$ws.Range("A17").Value = "Last Name"
$ws.Range("B17").Value = "Last"
$ws.Range("B19").Value = "Address1"
$ws.Range("B20").Value = "Address2"
$ws.Range("A21").Value = "City"
$ws.Range("A22").Value = "State"
$ws.Range("B22").Value = "Auckland"
$ws.Range("B23").Value = 335443
$ws.Range("A24").Value = "Country"
$ws.Range("B25").Value = "Yes"
$ws.Range("B26").Value = "Address1"
$ws.Range("B27").Value = "Address2"
$ws.Range("A28").Value = "City"
$ws.Range("B28").Value = "Auckland"
$ws.Range("A29").Value = "State"
$ws.Range("B29").Value = "Auckland"
$ws.Range("A30").Value = "Postal Code"
$ws.Range("B30").Value = 335443
$ws.Range("A31").Value = "Country"
$ws.Range("B31").Value = "NEW ZEALAND"

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- View state: scroll window + select final cell ---
$ws.Activate()
$ws.Range("A32").Select()
$excel.ActiveWindow.ScrollRow = 28
